$d = $word.ActiveDocument

# 1. Title / heading (appears twice: Heading1 at top, and bold run near bottom)
$d.Content.Find.Execute("Play Book of the Sphinx for Free – Review & Bonus Offers", $true, $false, $false, $false, $false, $true, 1, $false, "Play Book of the Sphinx Free | Exciting Egyptian Themed Slot Game", 2, 2)

# 2. "What we like" bullet list
$d.Content.Find.Execute("Engaging gameplay, suitable for different types of players", $true, $false, $false, $false, $false, $true, 1, $false, "Engaging gameplay with a simple interface", 2)
$d.Content.Find.Execute("Book of the Sphinx symbol, high payout of up to 4000 times initial bet", $true, $false, $false, $false, $false, $true, 1, $false, "High-quality graphics and smooth animations", 2)
$d.Content.Find.Execute("Thematic symbols and simple graphics add to the game's appeal", $true, $false, $false, $false, $false, $true, 1, $false, "Exciting Gamble feature for adrenaline junkies", 2)
$d.Content.Find.Execute("Gamble feature adds a layer of excitement to the game", $true, $false, $false, $false, $false, $true, 1, $false, "Suitable for players of all budgets", 2)

# 3. "What we don't like" bullet list
$d.Content.Find.Execute("Graphics may be too basic for some players", $true, $false, $false, $false, $false, $true, 1, $false, "Basic graphics may not appeal to everyone", 2)
$d.Content.Find.Execute("Maximum bet of €200 may not be high enough for high rollers", $true, $false, $false, $false, $false, $true, 1, $false, "Gamble feature can result in total loss of winnings", 2)

# 4. Meta description (italic run at the very end)
$d.Content.Find.Execute("Read our review of Book of the Sphinx online slot game. Play for free and discover bonus offers. Suitable for all players. Try your luck now!", $true, $false, $false, $false, $false, $true, 1, $false, "Experience the adventure of Ancient Egypt with Book of the Sphinx. Play for free and win big!", 2)
